$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species-record data between row 22 and row 23 for the
# columns that actually differ between the two records (the shared
# metadata columns - D, I, P, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY -
# stay untouched since they are identical for both rows).

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "S", "Z", "AB")

foreach ($col in $cols) {
    $addr22 = "$col`22"
    $addr23 = "$col`23"

    $val22 = $ws.Range($addr22).Value2
    $val23 = $ws.Range($addr23).Value2

    $ws.Range($addr22).Value2 = $val23
    $ws.Range($addr23).Value2 = $val22
}
